# Update the "Förändrad" date (column C) for rows 2-12 from
# serial date 45233 (2023-11-03) to 45243 (2023-11-13).
# Cell formatting/style is left untouched by only updating the value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 12; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $cell.Value = 45243
}
